$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("H2").Value = 2.7
$ws.Range("I2").Value = 4
$ws.Range("J2").Value = 1.2
$ws.Range("K2").Value = 4.33
$ws.Range("P2").Value = 1.9
$ws.Range("Q2").Value = 1.9
$ws.Range("R2").Value = 2.75
$ws.Range("S2").Value = 1.4
$ws.Range("X2").Value = 29
$ws.Range("Z2").Value = 4
$ws.Range("AA2").Value = 6
$ws.Range("AE2").Value = 6.5
$ws.Range("AJ2").Value = 67
$ws.Range("G3").Value = 2.2
$ws.Range("I3").Value = 3.5
$ws.Range("K3").Value = 7.5
$ws.Range("T3").Value = 6.5
$ws.Range("U3").Value = 9.5
$ws.Range("W3").Value = 21
$ws.Range("X3").Value = 21
$ws.Range("AC3").Value = 51
$ws.Range("AD3").Value = 401
$ws.Range("AE3").Value = 8.5
$ws.Range("G4").Value = 2.88
$ws.Range("I4").Value = 2.6
$ws.Range("T4").Value = 7.5
$ws.Range("U4").Value = 13
$ws.Range("AF4").Value = 11
$ws.Range("AH4").Value = 26
$ws.Range("AI4").Value = 23
$ws.Range("G8").Value = 1.95
$ws.Range("I8").Value = 4.5
$ws.Range("J8").Value = 1.11
$ws.Range("K8").Value = 6.5
$ws.Range("P8").Value = 1.53
$ws.Range("Q8").Value = 2.38
$ws.Range("U8").Value = 8
$ws.Range("W8").Value = 17
$ws.Range("Z8").Value = 6.5
$ws.Range("AB8").Value = 19
$ws.Range("AC8").Value = 67
$ws.Range("AG8").Value = 15
$ws.Range("U9").Value = 11
$ws.Range("AC9").Value = 67
$ws.Range("AE9").Value = 7.5
$ws.Range("G10").Value = 1.9
$ws.Range("H10").Value = 3.2
$ws.Range("J10").Value = 1.08
$ws.Range("K10").Value = 8
$ws.Range("L10").Value = 1.4
$ws.Range("M10").Value = 2.75
$ws.Range("N10").Value = 2.3
$ws.Range("O10").Value = 1.6
$ws.Range("P10").Value = 1.5
$ws.Range("Q10").Value = 2.5
$ws.Range("T10").Value = 6
$ws.Range("V10").Value = 9
$ws.Range("X10").Value = 17
$ws.Range("AA10").Value = 6.5
$ws.Range("AD10").Value = 1250
$ws.Range("AG10").Value = 15
$ws.Range("H11").Value = 3.2
$ws.Range("Z11").Value = 7.5
$ws.Range("AB11").Value = 19
$ws.Range("G12").Value = 1.45
$ws.Range("H12").Value = 4
$ws.Range("I12").Value = 6.5
$ws.Range("N12").Value = 1.85
$ws.Range("O12").Value = 1.95
$ws.Range("R12").Value = 2
$ws.Range("S12").Value = 1.73
$ws.Range("T12").Value = 6.5
$ws.Range("U12").Value = 6.5
$ws.Range("AC12").Value = 67
$ws.Range("O14").Value = 1.8
$ws.Range("P14").Value = 1.44
$ws.Range("Q14").Value = 2.63
$ws.Range("R14").Value = 1.8
$ws.Range("S14").Value = 1.91
$ws.Range("T14").Value = 8.5
$ws.Range("AB14").Value = 15
$ws.Range("AC14").Value = 51
$ws.Range("AD14").Value = 251
$ws.Range("AE14").Value = 8.5
$ws.Range("G15").Value = 3
$ws.Range("I15").Value = 2.25
$ws.Range("N15").Value = 1.93
$ws.Range("O15").Value = 1.93
$ws.Range("T15").Value = 10
$ws.Range("W15").Value = 34
$ws.Range("AF15").Value = 11
$ws.Range("AH15").Value = 21
$ws.Range("AJ15").Value = 26
$ws.Range("G16").Value = 19.5
$ws.Range("I16").Value = 1.12
$ws.Range("N16").Value = 1.41
$ws.Range("O16").Value = 2.7
$ws.Range("R16").Value = 2.38
$ws.Range("S16").Value = 1.52
$ws.Range("T16").Value = 45
$ws.Range("U16").Value = 200
$ws.Range("V16").Value = 55
$ws.Range("X16").Value = 300
$ws.Range("Z16").Value = 16
$ws.Range("AA16").Value = 12.5
$ws.Range("AB16").Value = 27
$ws.Range("AE16").Value = 7.1
$ws.Range("AH16").Value = 5.3
$ws.Range("AJ16").Value = 28
$ws.Range("G17").Value = 2.18
$ws.Range("H17").Value = 3.2
$ws.Range("I17").Value = 2.95
$ws.Range("L17").Value = 1.36
$ws.Range("M17").Value = 2.9
$ws.Range("N17").Value = 2.02
$ws.Range("P17").Value = 1.39
$ws.Range("Q17").Value = 2.45
$ws.Range("R17").Value = 1.88
$ws.Range("S17").Value = 1.83
$ws.Range("T17").Value = 6
$ws.Range("W17").Value = 16.5
$ws.Range("Y17").Value = 25
$ws.Range("Z17").Value = 8.25
$ws.Range("AA17").Value = 5.4
$ws.Range("AC17").Value = 60
$ws.Range("AD17").Value = 450
$ws.Range("AE17").Value = 7
$ws.Range("AF17").Value = 11.75
$ws.Range("AH17").Value = 28
$ws.Range("AJ17").Value = 30
$ws.Range("G18").Value = 1.91
$ws.Range("I18").Value = 4.33
$ws.Range("J18").Value = 1.06
$ws.Range("K18").Value = 10
$ws.Range("N18").Value = 2.08
$ws.Range("O18").Value = 1.73
$ws.Range("R18").Value = 1.95
$ws.Range("S18").Value = 1.8
$ws.Range("U18").Value = 8.5
$ws.Range("AD18").Value = 351
$ws.Range("AF18").Value = 21
$ws.Range("G19").Value = 1.45
$ws.Range("I19").Value = 6.5
$ws.Range("L19").Value = 1.25
$ws.Range("M19").Value = 3.75
$ws.Range("N19").Value = 1.75
$ws.Range("O19").Value = 2.05
$ws.Range("W19").Value = 10
$ws.Range("AA19").Value = 8.5
$ws.Range("AB19").Value = 19
$ws.Range("AC19").Value = 51
$ws.Range("AD19").Value = 301
$ws.Range("AH19").Value = 67
$ws.Range("AJ19").Value = 41
$ws.Range("L20").Value = 1.29
$ws.Range("M20").Value = 3.5
$ws.Range("N20").Value = 1.9
$ws.Range("O20").Value = 1.9
$ws.Range("G21").Value = 3.4
$ws.Range("I21").Value = 2.25
$ws.Range("W21").Value = 41
$ws.Range("X21").Value = 34
$ws.Range("AE21").Value = 6
$ws.Range("AF21").Value = 9.5
